$d = $word.ActiveDocument

# Locate "representative embedded system" in the introductory paragraph so we
# know exactly where the new text will be inserted (needed below to re-seat
# the _GoBack bookmark at the right spot).
$locate = $d.Content
$found = $locate.Find.Execute(
    "representative embedded system", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $insertPos   = $locate.End
    $bookmarkPos = $insertPos + 5   # length of " (an " -> right before "IoT device"

    # Expand "... representative embedded system." to
    # "... representative embedded system (an IoT device)."
    [void]$d.Content.Find.Execute(
        "representative embedded system", $true, $false, $false, $false, $false,
        $true, 1, $false, "representative embedded system (an IoT device)", 2)

    # Word re-seats the hidden _GoBack bookmark at the location of the most
    # recent edit. Reproduce that here: drop the old bookmark and re-add it
    # right after "an " / before "IoT device" inside the text we just typed.
    try {
        $goBack = $d.Bookmarks.Item("_GoBack")
        [void]$goBack.Delete()
    } catch {
    }
    [void]$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
}
